# Applies the "organized and ran data from 20250513" edit:
# Duplicates the existing 20250506 data block (rows 2-31) into rows 32-61,
# but stamped with the new run date 20250513. Then moves the active
# selection to C57 (the last-edited cell in that new block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 20250513

# Each entry: target row, sample (B), well (C), type (D), family (E, optional)
$rows = @(
    @{ R = 32; B = "A1";     C = "A01"; D = "sample"; E = "A" },
    @{ R = 33; B = "A2";     C = "A02"; D = "sample"; E = "A" },
    @{ R = 34; B = "A3";     C = "A03"; D = "sample"; E = "A" },
    @{ R = 35; B = "A4";     C = "A04"; D = "sample"; E = "A" },
    @{ R = 36; B = "A5";     C = "A05"; D = "sample"; E = "A" },
    @{ R = 37; B = "B1";     C = "B01"; D = "sample"; E = "B" },
    @{ R = 38; B = "B2";     C = "B02"; D = "sample"; E = "B" },
    @{ R = 39; B = "B3";     C = "B03"; D = "sample"; E = "B" },
    @{ R = 40; B = "B4";     C = "B04"; D = "sample"; E = "B" },
    @{ R = 41; B = "B5";     C = "B05"; D = "sample"; E = "B" },
    @{ R = 42; B = "C1";     C = "C01"; D = "sample"; E = "C" },
    @{ R = 43; B = "C2";     C = "C02"; D = "sample"; E = "C" },
    @{ R = 44; B = "C3";     C = "C03"; D = "sample"; E = "C" },
    @{ R = 45; B = "C4";     C = "C04"; D = "sample"; E = "C" },
    @{ R = 46; B = "C5";     C = "C05"; D = "sample"; E = "C" },
    @{ R = 47; B = "D1";     C = "D01"; D = "sample"; E = "D" },
    @{ R = 48; B = "D2";     C = "D02"; D = "sample"; E = "D" },
    @{ R = 49; B = "D3";     C = "D03"; D = "sample"; E = "D" },
    @{ R = 50; B = "D4";     C = "D04"; D = "sample"; E = "D" },
    @{ R = 51; B = "D5";     C = "D05"; D = "sample"; E = "D" },
    @{ R = 52; B = "E1";     C = "E01"; D = "sample"; E = "E" },
    @{ R = 53; B = "E2";     C = "E02"; D = "sample"; E = "E" },
    @{ R = 54; B = "E3";     C = "E03"; D = "sample"; E = "E" },
    @{ R = 55; B = "E4";     C = "E04"; D = "sample"; E = "E" },
    @{ R = 56; B = "E5";     C = "E05"; D = "sample"; E = "E" },
    @{ R = 57; B = "Blank1"; C = "H01"; D = "blank";  E = $null },
    @{ R = 58; B = "Blank2"; C = "H02"; D = "blank";  E = $null },
    @{ R = 59; B = "Blank3"; C = "H03"; D = "blank";  E = $null },
    @{ R = 60; B = "Blank4"; C = "H04"; D = "blank";  E = $null },
    @{ R = 61; B = "Blank5"; C = "H05"; D = "blank";  E = $null }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = $newDate
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    if ($row.E) {
        $ws.Range("E$r").Value = $row.E
    }
}

# Match the font color already used on the "sample"/"well" columns elsewhere
# in the sheet (black RGB rather than the theme color) so no new style is
# introduced.
$ws.Range("B32:C61").Font.Color = 0

# The author's last selection after entering this block.
$ws.Range("C57").Select() | Out-Null

Write-Host "Added 20250513 data block (rows 32-61) and updated selection."
